# TestTask - rename the "Type" sheet to "ProductType" and switch the
# active/selected tab from "Category" to the renamed sheet, matching the
# workbook.xml / sheet3.xml / sheet4.xml changes in the target diff.

$wb = $excel.ActiveWorkbook

# Rename sheet "Type" -> "ProductType" (sheets/workbook.xml <sheet name="..."/>)
$wsProductType = $wb.Worksheets.Item("Type")
$wsProductType.Name = "ProductType"

# Make "ProductType" the active sheet. This:
#   - updates <workbookView ... activeTab="3" .../> in workbook.xml
#   - moves tabSelected="1" off the "Category" sheetView (sheet3.xml)
#   - adds tabSelected="1" to the "ProductType" sheetView (sheet4.xml)
$wsProductType.Activate()
